$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 4475
$ws.Range("E2").Value = 48
$ws.Range("F2").Value = 48
$ws.Range("G2").Value = -26
$ws.Range("H2").Value = -46
$ws.Range("I2").Value = -58
$ws.Range("J2").Value = 12
$ws.Range("K2").Value = 5157
$ws.Range("L2").Value = 3210
$ws.Range("M2").Value = 1946
$ws.Range("N2").Value = 1704
$ws.Range("O2").Value = 242
$ws.Range("P2").Value = 346
$ws.Range("Q2").Value = -72
$ws.Range("R2").Value = -94
$ws.Range("S2").Value = -90
$ws.Range("T2").Value = 208
$ws.Range("U2").Value = -280
$ws.Range("V2").Value = 2628
$ws.Range("W2").Value = 1.08
$ws.Range("X2").Value = -1.03
$ws.Range("Y2").Value = -3.57
$ws.Range("Z2").Value = -0.88
$ws.Range("AA2").Value = 164.95
$ws.Range("AB2").Value = 415.28
$ws.Range("AC2").Value = -837
$ws.Range("AD2").Value = -21.39
$ws.Range("AE2").Value = 27880
$ws.Range("AF2").Value = 0.64
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 6976568

# Row 3
$ws.Range("D3").Value = 3963
$ws.Range("E3").Value = -126
$ws.Range("F3").Value = -126
$ws.Range("G3").Value = -526
$ws.Range("H3").Value = -564
$ws.Range("I3").Value = -476
$ws.Range("J3").Value = -88
$ws.Range("K3").Value = 4756
$ws.Range("L3").Value = 2946
$ws.Range("M3").Value = 1810
$ws.Range("N3").Value = 1276
$ws.Range("O3").Value = 534
$ws.Range("P3").Value = 346
$ws.Range("Q3").Value = -166
$ws.Range("R3").Value = -123
$ws.Range("S3").Value = 275
$ws.Range("T3").Value = 179
$ws.Range("U3").Value = -345
$ws.Range("V3").Value = 2465
$ws.Range("W3").Value = -3.18
$ws.Range("X3").Value = -14.22
$ws.Range("Y3").Value = -31.94
$ws.Range("Z3").Value = -11.37
$ws.Range("AA3").Value = 162.79
$ws.Range("AB3").Value = 273.9
$ws.Range("AC3").Value = -6822
$ws.Range("AD3").Value = -1.92
$ws.Range("AE3").Value = 20878
$ws.Range("AF3").Value = 0.63
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 6976568

# Row 4
$ws.Range("D4").Value = 4152
$ws.Range("E4").Value = 74
$ws.Range("F4").Value = 74
$ws.Range("G4").Value = -137
$ws.Range("H4").Value = -174
$ws.Range("I4").Value = -211
$ws.Range("J4").Value = 37
$ws.Range("K4").Value = 4605
$ws.Range("L4").Value = 3004
$ws.Range("M4").Value = 1601
$ws.Range("N4").Value = 1023
$ws.Range("O4").Value = 579
$ws.Range("P4").Value = 346
$ws.Range("Q4").Value = 238
$ws.Range("R4").Value = -81
$ws.Range("S4").Value = -127
$ws.Range("T4").Value = 73
$ws.Range("U4").Value = 165
$ws.Range("V4").Value = 2347
$ws.Range("W4").Value = 1.78
$ws.Range("X4").Value = -4.2
$ws.Range("Y4").Value = -18.39
$ws.Range("Z4").Value = -3.72
$ws.Range("AA4").Value = 187.56
$ws.Range("AB4").Value = 211.4
$ws.Range("AC4").Value = -3030
$ws.Range("AD4").Value = -3.11
$ws.Range("AE4").Value = 16733
$ws.Range("AF4").Value = 0.5600000000000001
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 6976568

# Row 5
$ws.Range("D5").Value = 1014
$ws.Range("E5").Value = -77
$ws.Range("F5").Value = -77
$ws.Range("G5").Value = -216
$ws.Range("H5").Value = -202
$ws.Range("I5").Value = -185
$ws.Range("J5").Value = -17
$ws.Range("K5").Value = 4099
$ws.Range("L5").Value = 2716
$ws.Range("M5").Value = 1383
$ws.Range("N5").Value = 623
$ws.Range("O5").Value = 759
$ws.Range("P5").Value = 346
$ws.Range("Q5").Value = 142
$ws.Range("R5").Value = -15
$ws.Range("S5").Value = -116
$ws.Range("T5").Value = 33
$ws.Range("U5").Value = 110
$ws.Range("V5").Value = 2219
$ws.Range("W5").Value = -7.64
$ws.Range("X5").Value = -19.91
$ws.Range("Y5").Value = -22.42
$ws.Range("Z5").Value = -4.64
$ws.Range("AA5").Value = 196.45
$ws.Range("AB5").Value = 159.34
$ws.Range("AC5").Value = -2645
$ws.Range("AD5").Value = -2.7
$ws.Range("AE5").Value = 10200
$ws.Range("AF5").Value = 0.7
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 6976568

# Row 6
$ws.Range("D6").Value = 773
$ws.Range("E6").Value = -478
$ws.Range("F6").Value = -478
$ws.Range("G6").Value = -717
$ws.Range("H6").Value = -209
$ws.Range("I6").Value = -181
$ws.Range("K6").Value = 1703
$ws.Range("L6").Value = 1356
$ws.Range("M6").Value = 347
$ws.Range("N6").Value = 493
$ws.Range("P6").Value = 349
$ws.Range("Q6").Value = -22
$ws.Range("R6").Value = 591
$ws.Range("S6").Value = -267
$ws.Range("T6").Value = 36
$ws.Range("U6").Value = -58
$ws.Range("V6").Value = 1029
$ws.Range("W6").Value = -61.84
$ws.Range("X6").Value = -27.08
$ws.Range("Y6").Value = -32.49
$ws.Range("Z6").Value = -7.22
$ws.Range("AA6").Value = 390.87
$ws.Range("AB6").Value = 106.99
$ws.Range("AC6").Value = -2600
$ws.Range("AD6").Value = -2.54
$ws.Range("AE6").Value = 8072
$ws.Range("AF6").Value = 0.82
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 6976569

# Rows 7-9: clear financial data columns, keep only A/B/C
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
